$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 29, which shifts the existing
# rows 29..146 down to 30..147 (preserving their formatting/styles).
$ws.Rows("29:29").Insert()

# Populate the newly inserted row 29 with the new price record.
$ws.Range("A29").Value = 5
$ws.Range("B29").Value = "Macroferia Regional de Talca"
$ws.Range("C29").Value = "Maule"
$ws.Range("D29").Value = 44991
$ws.Range("E29").Value = 7
$ws.Range("F29").Value = 100112001
$ws.Range("G29").Value = "Berenjena"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 200
$ws.Range("K29").Value = 8000
$ws.Range("L29").Value = 8000
$ws.Range("M29").Value = 8000
$ws.Range("N29").Value = "`$/caja 50 unidades"
$ws.Range("O29").Value = "Región del Maule"
$ws.Range("P29").Value = 160
$ws.Range("Q29").Value = 50
$ws.Range("R29").Value = "Hortaliza"
